$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column J (JuristPerson), shifting it to K.
$ws.Columns.Item(10).Insert()

# New header for the inserted column J.
$ws.Range("J1").Value = "Personensteuer"

# New data values for column J, rows 2-5.
$ws.Range("J2").Value = 10
$ws.Range("J3").Value = 10
$ws.Range("J4").Value = 20
$ws.Range("J5").Value = 20

# Update selection to match the target state.
$ws.Range("J11").Select()
